$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.882.56'
$ws.Range('E2').Value = '  -1.38%  '
$ws.Range('D3').Value = '2.684.66'
$ws.Range('E3').Value = '  -1.87%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '554.99'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  -1.47%  '
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '158.09'
$ws.Range('D6').Style = $origStyle
$ws.Range('E6').Value = '  -0.73%  '
$origStyle = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.589'
$ws.Range('D8').Style = $origStyle
$ws.Range('E8').Value = '  -1.08%  '
$ws.Range('E9').Value = '  -2.71%  '
$ws.Range('E10').Value = '  -2.21%  '
$ws.Range('E11').Value = '  -2.62%  '
$ws.Range('E12').Value = '  -3.73%  '
$ws.Range('D13').Value = '3.157.02'
$ws.Range('E13').Value = '  -1.92%  '
$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '26.62'
$ws.Range('D14').Style = $origStyle
$ws.Range('E14').Value = '  -1.01%  '
$ws.Range('D15').Value = '62.790.97'
$ws.Range('E15').Value = '  -1.29%  '
$ws.Range('E16').Value = '  -1.58%  '
$ws.Range('D17').Value = '2.682.26'
$ws.Range('E17').Value = '  -2.10%  '
$origStyle = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '11.89'
$ws.Range('D18').Style = $origStyle
$ws.Range('E18').Value = '  -3.68%  '
$origStyle = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.62'
$ws.Range('D19').Style = $origStyle
$ws.Range('E19').Value = '  -2.56%  '
$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '345.47'
$ws.Range('D20').Style = $origStyle
$ws.Range('E20').Value = '  -2.22%  '
$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.24'
$ws.Range('D21').Style = $origStyle
$ws.Range('E21').Value = '  -4.95%  '
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('E23').Value = '  -2.78%  '
$origStyle = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '63.17'
$ws.Range('D24').Style = $origStyle
$ws.Range('E24').Value = '  -1.50%  '
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('E26').Value = '  -0.01%  '
$origStyle = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.20'
$ws.Range('D27').Style = $origStyle
$ws.Range('E27').Value = '  -1.91%  '
$ws.Range('E28').Value = '  +9.08%  '
$ws.Range('D29').Value = '0.0₃0858'
$ws.Range('E29').Value = '  -5.05%  '
$origStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.27'
$ws.Range('D30').Style = $origStyle
$ws.Range('E30').Value = '  +0.91%  '
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '163.94'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  +0.39%  '
$ws.Range('E33').Value = '  +0.97%  '
$ws.Range('E34').Value = '  +0.49%  '
$ws.Range('E35').Value = '  +0.00%  '
$origStyle = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '19.49'
$ws.Range('D36').Style = $origStyle
$ws.Range('E36').Value = '  -2.71%  '
$ws.Range('E37').Value = '  +0.09%  '
$origStyle = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '349.83'
$ws.Range('D38').Style = $origStyle
$ws.Range('E38').Value = '  +1.69%  '
$origStyle = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.27'
$ws.Range('D39').Style = $origStyle
$ws.Range('E39').Value = '  +0.08%  '
$ws.Range('E40').Value = '  -3.25%  '
$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.99'
$ws.Range('D41').Style = $origStyle
$ws.Range('E41').Value = '  -2.26%  '
$origStyle = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '38.35'
$ws.Range('D42').Style = $origStyle
$ws.Range('E42').Value = '  -0.20%  '
$origStyle = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '20.87'
$ws.Range('D43').Style = $origStyle
$ws.Range('E43').Value = '  -4.19%  '
$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '20.17'
$ws.Range('D44').Style = $origStyle
$ws.Range('E44').Value = '  -3.98%  '
$ws.Range('E45').Value = '  -0.93%  '
$ws.Range('E46').Value = '  -3.64%  '
$ws.Range('E47').Value = '  -0.09%  '
$origStyle = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '10.99'
$ws.Range('D48').Style = $origStyle
$ws.Range('E48').Value = '  -0.55%  '
$origStyle = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0971'
$ws.Range('D49').Style = $origStyle
$ws.Range('E50').Value = '  -3.01%  '
$origStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '128.72'
$ws.Range('D51').Style = $origStyle
$ws.Range('E51').Value = '  -4.14%  '
